$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

# The commit inserts a new (blank) column before column N on the
# "Repayment Schedule" sheet. This shifts the old "In Advance"/"Late"
# header cells (N1/O1 area) and the per-row "In Advance"/"Outstanding"
# data one column to the right (old N -> O, old P -> Q), leaving the
# freshly inserted column N blank.
$ws.Columns("N").Insert()

# Match the saved selection state from the edited workbook.
$ws.Range("S11").Select()
